# "Versão final - Backlog e Menual de Instalação"
#
# Marks the remaining SP3 backlog items as delivered (STATUS: PENDENTE -> ENTREGUE)
# on the Backlog_AGGRAN sheet. All the downstream SUMIFS totals (M6/N6 on
# Backlog_AGGRAN, and the mirrored cells on "Lista - validações") recompute
# automatically from these input changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog_AGGRAN")
$ws.Activate()

# Rows whose STATUS column (H) flips from PENDENTE to ENTREGUE.
$rows = @(21, 22, 23, 24, 25, 27, 28, 31, 33, 34, 35, 36)
foreach ($r in $rows) {
    $ws.Range("H$r").Value = "ENTREGUE"
}

# Leave the sheet's selection where the author ended up after this edit.
$ws.Range("B10").Select()
